# GitHub Actions data refresh for the cryptos list: new prices / 1h-volume
# percentages for most rows, plus Polkadot and WrappedEther swapping ranks
# (rows 12 <-> 13: Coin name, Link and Price all move with the swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value. Price-column values that are ambiguous with plain
# numbers (a single '.' and otherwise only digits) carry a leading "'"
# so Excel stores them as text, exactly like the original cell content.
$updates = @(
    @('D2', '26.418.96'),
    @('E2', '  +0.44%  '),
    @('D3', '1.698.63'),
    @('E3', '  +1.06%  '),
    @('D4', '''1.009'),
    @('E4', '  +0.18%  '),
    @('D5', '''218.59'),
    @('E5', '  -0.08%  '),
    @('D6', '''0.5479'),
    @('E6', '  +3.94%  '),
    @('D9', '''0.06451'),
    @('E9', '  +0.21%  '),
    @('D10', '''22.00'),
    @('E10', '  -0.35%  '),
    @('E11', '  +2.31%  '),
    @('B12', 'Polkadot'),
    @('C12', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'),
    @('D12', '''4.561'),
    @('E12', '  +0.25%  '),
    @('B13', 'WrappedEther'),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('D13', '1.691.92'),
    @('E13', '  +0.09%  '),
    @('D14', '''0.5859'),
    @('E14', '  +0.79%  '),
    @('D15', '''0.000008410'),
    @('E15', '  -1.22%  '),
    @('D16', '''65.70'),
    @('E16', '  +1.83%  '),
    @('D17', '26.459.77'),
    @('E17', '  +0.43%  '),
    @('D18', '''4.950'),
    @('E18', '  +0.28%  '),
    @('D19', '''1.009'),
    @('E19', '  +0.19%  '),
    @('E20', '  +1.00%  '),
    @('D21', '''191.63'),
    @('D22', '''6.273'),
    @('E22', '  +0.91%  '),
    @('D23', '''1.009'),
    @('E23', '  +0.11%  '),
    @('D24', '''148.66'),
    @('E24', '  +2.34%  '),
    @('D25', '''0.1310'),
    @('E25', '  +5.01%  '),
    @('D26', '''7.933'),
    @('E26', '  +1.99%  '),
    @('D27', '''15.80'),
    @('E27', '  -0.24%  '),
    @('D28', '''0.06231'),
    @('E28', '  -6.20%  '),
    @('D29', '''1.392'),
    @('E29', '  +2.49%  '),
    @('E30', '  -0.08%  '),
    @('D31', '''3.614'),
    @('E31', '  +1.01%  '),
    @('D32', '''3.597'),
    @('E32', '  +0.11%  '),
    @('D33', '''1.688'),
    @('E33', '  +1.53%  '),
    @('D34', '''1.038'),
    @('E34', '  +1.00%  '),
    @('D35', '''0.6162'),
    @('E35', '  -0.90%  '),
    @('D36', '''2.410'),
    @('E36', '  +0.57%  '),
    @('D37', '''2.760'),
    @('E37', '  +1.23%  '),
    @('D38', '''0.01661'),
    @('E38', '  +2.27%  '),
    @('D39', '1.118.47'),
    @('D40', '''6.121'),
    @('E40', '  -4.99%  '),
    @('D41', '''0.8821'),
    @('E41', '  +0.26%  '),
    @('D42', '''1.015'),
    @('E42', '  +0.04%  '),
    @('D43', '''101.21'),
    @('E43', '  +0.48%  '),
    @('D44', '1.847.99'),
    @('E44', '  +0.87%  '),
    @('E45', '  -4.28%  '),
    @('D46', '''57.69'),
    @('E46', '  +1.21%  '),
    @('D47', '''8.235'),
    @('E47', '  +1.19%  '),
    @('D48', '''1.009'),
    @('E48', '  -0.12%  '),
    @('E49', '  +0.15%  '),
    @('D50', '''6.132'),
    @('E50', '  +0.95%  '),
    @('D51', '''0.4300'),
    @('E51', '  +0.00%  '),
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $range = $ws.Range($cellRef)
    $wasForcedText = $newValue.StartsWith("'")
    $range.Value = $newValue
    if ($wasForcedText) {
        # Drop the quote-prefix style Excel applies to text-forced numeric
        # entries so the cell's formatting matches the original (unstyled) cell.
        $range.Style = 'Normal'
    }
}
